# Excel COM-interop edit script
# Updates the "average" sheet's consensus/dense comparison table:
#  - refreshes all metric rows (new values, new rows, reordered metrics)
#  - highlights the "dense" comparison rows (ibes_1_fwdepsqcut /
#    ibes_1_fwdepsqcut_dense) in bold/purple
#  - re-applies the column-C AutoFilter (hides rows whose `len` isn't
#    one of the "full sample" counts) and widens column A to fit the
#    longest index label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Write header + data rows ----
$ws.Cells.Item(1,1).Value = "index"
$ws.Cells.Item(1,2).Value = "consensus"
$ws.Cells.Item(1,3).Value = "len"
$ws.Cells.Item(1,4).Value = "lgbm"

$ws.Cells.Item(2,1).Value = "ibes_6_ni"
$ws.Cells.Item(2,2).Value = 0.007477441611400649
$ws.Cells.Item(2,3).Value = 568
$ws.Cells.Item(2,4).Value = 0.006952197441597438

$ws.Cells.Item(3,1).Value = "ibes_2_ni"
$ws.Cells.Item(3,2).Value = 0.008604109548964898
$ws.Cells.Item(3,3).Value = 14467
$ws.Cells.Item(3,4).Value = 0.007461390831007631

$ws.Cells.Item(4,1).Value = "ibes_2_fwdeps"
$ws.Cells.Item(4,2).Value = 0.008604109548964898
$ws.Cells.Item(4,3).Value = 14467
$ws.Cells.Item(4,4).Value = 0.007470316543135013

$ws.Cells.Item(5,1).Value = "ibes_2_fwdepsqcut"
$ws.Cells.Item(5,2).Value = 0.008604109548964898
$ws.Cells.Item(5,3).Value = 14467
$ws.Cells.Item(5,4).Value = 0.008215291392886251

$ws.Cells.Item(6,1).Value = "ibes_1_fwdepsqcut"
$ws.Cells.Item(6,2).Value = 0.008610171467262949
$ws.Cells.Item(6,3).Value = 13776
$ws.Cells.Item(6,4).Value = 0.008239669139180731

$ws.Cells.Item(7,1).Value = "ibes_6_niqcut"
$ws.Cells.Item(7,2).Value = 0.008673141980908865
$ws.Cells.Item(7,3).Value = 13402
$ws.Cells.Item(7,4).Value = 0.008290320575140078

$ws.Cells.Item(8,1).Value = "ibes_1_fwdepsqcut_dense"
$ws.Cells.Item(8,2).Value = 0.008658944184046573
$ws.Cells.Item(8,3).Value = 11650
$ws.Cells.Item(8,4).Value = 0.009065711096231268

$ws.Cells.Item(9,1).Value = "ni_2_niqcut"
$ws.Cells.Item(9,2).Value = 0.008528350971131086
$ws.Cells.Item(9,3).Value = 2695
$ws.Cells.Item(9,4).Value = 0.01262906986070185

$ws.Cells.Item(10,1).Value = "ni_6_fwdeps"
$ws.Cells.Item(10,2).Value = 0.009038742197873241
$ws.Cells.Item(10,3).Value = 4248
$ws.Cells.Item(10,4).Value = 0.0126430822285556

$ws.Cells.Item(11,1).Value = "ni_1_epsqcut"
$ws.Cells.Item(11,2).Value = 0.00930690016118596
$ws.Cells.Item(11,3).Value = 14467
$ws.Cells.Item(11,4).Value = 0.01322805762894923

$ws.Cells.Item(12,1).Value = "ni_2_epsqcut"
$ws.Cells.Item(12,2).Value = 0.00938408729045367
$ws.Cells.Item(12,3).Value = 14289
$ws.Cells.Item(12,4).Value = 0.01340126360609522

$ws.Cells.Item(13,1).Value = "ni_2_fwdeps"
$ws.Cells.Item(13,2).Value = 0.008644807991816253
$ws.Cells.Item(13,3).Value = 8817
$ws.Cells.Item(13,4).Value = 0.01341510294821813

$ws.Cells.Item(14,1).Value = "ni_1_fwdepsqcut"
$ws.Cells.Item(14,2).Value = 0.00930690016118596
$ws.Cells.Item(14,3).Value = 14467
$ws.Cells.Item(14,4).Value = 0.01367246194968542

$ws.Cells.Item(15,1).Value = "ni_2_fwdepsqcut"
$ws.Cells.Item(15,2).Value = 0.00938408729045367
$ws.Cells.Item(15,3).Value = 14289
$ws.Cells.Item(15,4).Value = 0.01373872658434699

$ws.Cells.Item(16,1).Value = "ni_6_epsqcut"
$ws.Cells.Item(16,2).Value = 0.009486213628492479
$ws.Cells.Item(16,3).Value = 14246
$ws.Cells.Item(16,4).Value = 0.01400753774994666

$ws.Cells.Item(17,1).Value = "ni_6_fwdepsqcut"
$ws.Cells.Item(17,2).Value = 0.009492788784155153
$ws.Cells.Item(17,3).Value = 14279
$ws.Cells.Item(17,4).Value = 0.01403873064787338

$ws.Cells.Item(18,1).Value = "ni_6_niqcut"
$ws.Cells.Item(18,2).Value = 0.01228465704399278
$ws.Cells.Item(18,3).Value = 172
$ws.Cells.Item(18,4).Value = 0.01536086770476644

$ws.Cells.Item(19,1).Value = "ibes_2_niqcut"
$ws.Cells.Item(19,2).Value = 0.00805737483750827
$ws.Cells.Item(19,3).Value = 9416
$ws.Cells.Item(19,4).Value = 0.1837138298836869


# ---- Highlight the "dense" comparison rows (6 and 8) in bold/regular purple ----
$purple = 10498160  # RGB(0x70,0x30,0xA0) -> FF7030A0

$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").Font.Color = $purple
$ws.Range("B6:D6").Font.Bold = $false
$ws.Range("B6:D6").Font.Color = $purple

$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").Font.Color = $purple
$ws.Range("B8:D8").Font.Bold = $false
$ws.Range("B8:D8").Font.Color = $purple

# ---- Widen column A to fit the longest index label ----
$ws.Columns.Item(1).ColumnWidth = 21

# ---- Re-apply the AutoFilter on column C (len), keep only the "full"
#      sample sizes -> this also hides the non-matching rows ----
$ws.Range("A1:D19").AutoFilter(3, @("11650","13402","13776","14246","14279","14289","14467"), 7)

# ---- Excel records the active AutoFilter range as a hidden, sheet
#      scoped defined name ----
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=average!`$A`$1:`$D`$19")
$fdb.Visible = $false

# ---- Restore the view state (zoom + selected cell) ----
$ws.Application.ActiveWindow.Zoom = 249
$ws.Range("C5").Select()

Write-Host "edit complete"
